$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 26 ("RM 232") entirely - remaining rows shift up.
$ws.Rows.Item(26).Delete()
# After that deletion, the row that was "SC 92" is now row 27 - delete it too.
$ws.Rows.Item(27).Delete()

# Apply the remaining individual cell edits (final row numbering, after the two deletions).
$ws.Range("C2").Value = 14.9
$ws.Range("C3").Value = ""
$ws.Range("C4").Value = ""
$ws.Range("D5").Value = ""
$ws.Range("D8").Value = -13.9
$ws.Range("D10").Value = -14.7
$ws.Range("C11").Value = 11.4
$ws.Range("D12").Value = ""
$ws.Range("C13").Value = ""
$ws.Range("D15").Value = -15.2
$ws.Range("D18").Value = ""
$ws.Range("D19").Value = ""
$ws.Range("C21").Value = 12.7
$ws.Range("C25").Value = ""
$ws.Range("D25").Value = -15.5
$ws.Range("D27").Value = -14.6
$ws.Range("B29").Value = ""
$ws.Range("D29").Value = ""
$ws.Range("B33").Value = -19.5
$ws.Range("C33").Value = 10.4
$ws.Range("D33").Value = ""
